$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"

$msg38033bfe = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/313a6df70a547950e789c36f699af4fa94cc03a7/e2e/38033bfe-6d7f-4df0-a5c9-e64f359b4f38.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/85b7b6df25521388087d9ab58db97b6ef370d37a/e2e/38033bfe-6d7f-4df0-a5c9-e64f359b4f38.md."
$msg741cf1f0 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/313a6df70a547950e789c36f699af4fa94cc03a7/e2e/741cf1f0-380b-4407-aacf-b3967867eddf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/85b7b6df25521388087d9ab58db97b6ef370d37a/e2e/741cf1f0-380b-4407-aacf-b3967867eddf.md."

# ---- Overview sheet: rows 4 & 5 switch to "Ready for handoff" with a new timestamp ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(4, 5).Value = $statusReady
$wsOverview.Cells.Item(4, 6).Value = $statusReady
$wsOverview.Cells.Item(4, 7).Value = "2016-08-17 04:23:58"
$wsOverview.Cells.Item(5, 5).Value = $statusReady
$wsOverview.Cells.Item(5, 6).Value = $statusReady
$wsOverview.Cells.Item(5, 7).Value = "2016-08-17 04:23:58"

# ---- zh-cn sheet: rows 4 & 5 Status/Datetime/Error Detail ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(4, 3).Value = $statusReady
$wsZhCn.Cells.Item(4, 8).Value = "2016-08-17 04:23:54"
$wsZhCn.Cells.Item(4, 16).Value = $msg38033bfe
$wsZhCn.Cells.Item(5, 3).Value = $statusReady
$wsZhCn.Cells.Item(5, 8).Value = "2016-08-17 04:23:54"
$wsZhCn.Cells.Item(5, 16).Value = $msg741cf1f0
# Error Detail column (P) widens to fit the new long messages
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---- de-de sheet: rows 4 & 5 Status/Datetime/Error Detail ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(4, 3).Value = $statusReady
$wsDeDe.Cells.Item(4, 8).Value = "2016-08-17 04:23:58"
$wsDeDe.Cells.Item(4, 16).Value = $msg38033bfe
$wsDeDe.Cells.Item(5, 3).Value = $statusReady
$wsDeDe.Cells.Item(5, 8).Value = "2016-08-17 04:23:58"
$wsDeDe.Cells.Item(5, 16).Value = $msg741cf1f0
# Error Detail column (P) widens to fit the new long messages
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
